$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "FAILED: $find"
    }
}

# Title: merge "Marker " + "WebApp" + " Overall Model" (drop proofErr wrapping) into one run.
Replace-Text "Marker WebApp Overall Model" "Marker WebApp Overall Model"

# Bullet: "Login page for the marker" -> "A login page for the marker"
Replace-Text "Login page for the marker" "A login page for the marker"

# Bullet: Marker have limited access...
Replace-Text "Marker have limited access to the feature of the web application." "Markers have limited access to features of the web application."

# Who: Marker is a member...
Replace-Text "Marker is a member of the organizer team. A desired user for the Web Application. Usually a teacher or a person involved with AMA or MATHEX." "A marker is a member of the organizer team. A desired user for the Web Application. Usually a teacher or a person involved with AMA or MATHEX. "

# What: Markers are responsible...
Replace-Text "Markers are responsible for marking the answer of each team that they are assigned to. Each marker, usually, marks two teams.  " "Markers are responsible for marking the answer of each team that they are assigned to. Each marker usually marks two teams.  "

# How: big rewrite
Replace-Text "Marker hands the question to a student. Then, student returns with the answer. If answer is correct, marker will note it down and compute the new total score. Marker hands to the student a new question to answer. This process continues until there is not more questions or the time is out. In case answer is wrong, student has the choice to try again or pass, if student pass a new question is given and marker notes it down as a pass." "The marker hands the question to a team who eventually return with an answer. If the answer is correct, the marker will mark it as ‘correct’. The marker will then hand the student the next question to answer. This process continues until there are no more questions or the competition time runs out. In case an answer is wrong the team has the choice to attempt am answer again or pass. If the team passes, the next question is given and the marker notes it down as a pass."

# Sequence diagram paragraph text update
Replace-Text "shows the required actions from the marker all the way to the database. It illustrates the login process, the selection of teams to mark, the process of chosen a selected team to mark and the possible interactions from the marker to with the webpage. " "shows the required actions from the marker all the way to the database. It illustrates the login process, the selection of teams to mark, the process of choosing a selected team to mark and the possible interactions from the marker with the webpage. "

# State diagram: "shows al states" -> "shows all states"
Replace-Text " – shows al states the webpage will hold as a result of each action performed by the marker." " – shows all states the webpage will hold as a result of each action performed by the marker."

# Heading: Marker interaction with webApp walk through (merge runs)
Replace-Text "Marker interaction with webApp walk through" "Marker interaction with webApp walk through"

# Markers are standard users...
Replace-Text "Markers are standard users of the application until they have logged in. In order to login, marker must access the login page. The login page can be reached indirectly by accessing the corresponding link. Or, directly through hyperlinks in the main page." "Markers are standard users of the application until they have logged in. In order to login, markers must access the login page. The login page can be reached by selecting ‘privileged user’ as the type of user. They must enter in the correct credentials to login."

# When logged-in...
Replace-Text "When logged-in. marker should be presented with a marker’s home page where they can select teams they want to mark. The number of teams to mark are usually 1 to 2, but it should be allowed more if required. " "When logged-in, markers should be presented with a marker’s home page where they can select teams they want to mark. The number of teams to mark are usually 1 to 2, but should be allowed more if required. "

# Each team will have its page...
Replace-Text "There will be thee action available in each page: " "There will be three actions available in each page: "

# Markers can choose teams...
Replace-Text "Markers can choose teams they want to mark ( alternatively, admin can pre-set them)" "Markers can choose teams they want to mark (alternatively, the admin can pre-set them)"

Write-Output "done"
